# Apply updated Pr(>F) p-values to column F following a rerun of the
# PERMANOVA permutation calculations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = 0.5987
    "F3"  = 0.0013
    "F4"  = 0.6901
    "F7"  = 0.3488
    "F8"  = 0.8515
    "F9"  = 0.6071
    "F12" = 0.5402
    "F13" = 0.9466
    "F14" = 0.2277
    "F17" = 0.7811
    "F18" = 0.3385
    "F19" = 0.8721
    "F22" = 0.1494
    "F23" = 0.0953
    "F24" = 0.2328
    "F27" = 0.2536
    "F28" = 0.6467
    "F29" = 0.7016
    "F32" = 0.8795
    "F33" = 0.0002
    "F34" = 0.9859
    "F37" = 0.3904
    "F38" = 0.0479
    "F39" = 0.4533
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
